# "Improved thesis road map" - nudge a batch of boxes/connectors in the
# roadmap diagram on slide 1 slightly to the left/right (only the X
# offset changes; Y/width/height are untouched).
#
# NOTE: the interop engine stores Shape.Left/.Top as a 32-bit float
# (in points) and truncates when converting back to EMU, so a plain
# "$emu / 12700" assignment can land 1 EMU short. Adding/subtracting
# half an EMU (in points) before the division compensates for that
# truncation and reproduces the exact target EMU value.
$EMU_PER_POINT = 12700

function Set-ShapeLeftEmu {
    param($shape, [double]$emu)
    $bump = 0.5
    if ($emu -lt 0) { $bump = -0.5 }
    $shape.Left = ($emu + $bump) / $EMU_PER_POINT
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

Set-ShapeLeftEmu $s.Shapes.Item("TextBox 5") 307627
Set-ShapeLeftEmu $s.Shapes.Item("TextBox 6") 1585189
Set-ShapeLeftEmu $s.Shapes.Item("Rectangle 16") 1620556
Set-ShapeLeftEmu $s.Shapes.Item("Rectangle 18") 352951
Set-ShapeLeftEmu $s.Shapes.Item("TextBox 20") -99231
Set-ShapeLeftEmu $s.Shapes.Item("Straight Connector 25") 1170641
Set-ShapeLeftEmu $s.Shapes.Item("Straight Connector 26") 1170641
Set-ShapeLeftEmu $s.Shapes.Item("Straight Connector 27") 2026004
Set-ShapeLeftEmu $s.Shapes.Item("TextBox 49") 352038
Set-ShapeLeftEmu $s.Shapes.Item("TextBox 50") 1596966
Set-ShapeLeftEmu $s.Shapes.Item("TextBox 590") 510859
Set-ShapeLeftEmu $s.Shapes.Item("Rectangle 593") 343800
Set-ShapeLeftEmu $s.Shapes.Item("Rectangle 594") 352951
Set-ShapeLeftEmu $s.Shapes.Item("Rectangle 595") 343800
Set-ShapeLeftEmu $s.Shapes.Item("Rectangle 596") 352951
Set-ShapeLeftEmu $s.Shapes.Item("TextBox 597") 486721
Set-ShapeLeftEmu $s.Shapes.Item("Straight Connector 599") 335562
Set-ShapeLeftEmu $s.Shapes.Item("Rectangle 618") 1637032
Set-ShapeLeftEmu $s.Shapes.Item("Rectangle 619") 1620556
Set-ShapeLeftEmu $s.Shapes.Item("Rectangle 621") 1637032
Set-ShapeLeftEmu $s.Shapes.Item("Rectangle 622") 1620556
Set-ShapeLeftEmu $s.Shapes.Item("TextBox 623") 1748063
Set-ShapeLeftEmu $s.Shapes.Item("Straight Connector 624") 1620556
Set-ShapeLeftEmu $s.Shapes.Item("TextBox 625") 1793970
Set-ShapeLeftEmu $s.Shapes.Item("Rectangle 626") 2904636
Set-ShapeLeftEmu $s.Shapes.Item("Rectangle 628") 2904636
